$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Update the J column (卡牌费用 / card fee) values to 0 for the filtered cards
$ws.Range("J8").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("J23").Value = 0

# Update the active selection on the sheet
$ws.Range("J26").Select()
